# Fruta / hortaliza, semanal
# Update "Fecha" (D) and "Volumen" (M) columns for rows 3-10 to reflect
# the re-shuffled weekly data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44330
$ws.Range("D4").Value = 44316
$ws.Range("D5").Value = 44327
$ws.Range("M5").Value = 60
$ws.Range("D6").Value = 44313
$ws.Range("D7").Value = 44302
$ws.Range("D8").Value = 44306
$ws.Range("M8").Value = 80
$ws.Range("D9").Value = 44322
$ws.Range("M9").Value = 60
$ws.Range("D10").Value = 44323
$ws.Range("M10").Value = 80
